# Roster shared by every monthly attendance sheet in this workbook.
$names = @(
    'AADYASHRI GUPTA',
    'AARADHYA JAIN',
    'AARNA MATHUR',
    'AATMIKA JAIN',
    'ANANYA MAHIRCHANDANI',
    'ANAYA VYAS',
    'ANVI JAIN',
    'ARADHYA GOYAL',
    'ARADHYA GURJAR',
    'AROHI KHANDELWAL',
    'BHOOMI KHATRI',
    'DIVIJA JOSHI',
    'DIVYANSHI',
    'DIVYANSHI JAIN',
    'GAURANSHI KHURANA',
    'GAURI KAUSHIK',
    'GUNJAN JOSHI',
    'HARSHITA GUPTA',
    'HIRAL',
    'HITANSHI JOGCHAND',
    'HIYA GOYAL',
    'KHUSHI JHA',
    'KRATI JAIN',
    'MANYA KHANDELWAL',
    'MUKTA KHATRI',
    'NAVYA GUPTA',
    'NEHAL JAIN',
    'PAAVANI SHARMA',
    'PRANJAL MAHESHWARI',
    'RUHEEN KHAN',
    'SAMIKSHA KHANDELWAL',
    'SANVI JAIN',
    'SATVIKA KANKARWAL',
    'SHIVI SETHI',
    'SHRASTI KHANGAROT',
    'SIDDHI JAIN',
    'SONAKSHI CHANDEL',
    'TANVI HIRANI',
    'TASHI JAIN',
    'VANSHIKA MAHAWAR',
    'YAJURVI KRISHNATREY',
    'YASHIKA KHANDPA'
)

$wb = $excel.ActiveWorkbook

# The workbook is rolling forward: drop the two oldest month sheets...
$wb.Worksheets("January").Delete() | Out-Null
$wb.Worksheets("February").Delete() | Out-Null

# ...clear out the stray "Jun 24" attendance column that had been filled in on
# the June sheet (P/A/H marks), so it goes back to a plain roster like the rest...
$june = $wb.Worksheets("June")
$june.Range("Y1:Y43").Clear() | Out-Null

# ...and add a new "November" sheet after October, with the same roster as the
# other months, so the tab list keeps pace with the calendar.
$october = $wb.Worksheets("October")
$november = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $october)
$november.Name = "November"
$november.Range("A1").Value = "NAME"
$november.Range("A1").Font.Bold = $true
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $november.Cells.Item($row, 1).Value = $names[$i]
}

# Make March (now the first sheet) the active tab, since January/February are gone.
$wb.Worksheets("March").Activate()
